# Update "想去人数" (want-to-go count) figures in column F across sheets.
# Mirrors the upstream gh-pages data refresh commit (output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1881
$ws1.Range("F3").Value  = 1511
$ws1.Range("F4").Value  = 879
$ws1.Range("F5").Value  = 765
$ws1.Range("F6").Value  = 13315
$ws1.Range("F7").Value  = 13181
$ws1.Range("F8").Value  = 1017
$ws1.Range("F10").Value = 21
$ws1.Range("F13").Value = 671
$ws1.Range("F14").Value = 2093
$ws1.Range("F15").Value = 64
$ws1.Range("F17").Value = 73
$ws1.Range("F20").Value = 246
$ws1.Range("F21").Value = 287
$ws1.Range("F22").Value = 419
$ws1.Range("F23").Value = 753
$ws1.Range("F24").Value = 15

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 34
$ws2.Range("F4").Value = 19

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 37

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1881
$ws4.Range("F4").Value  = 1511
$ws4.Range("F5").Value  = 879
$ws4.Range("F7").Value  = 765
$ws4.Range("F8").Value  = 13315
$ws4.Range("F9").Value  = 13181
$ws4.Range("F10").Value = 1017
$ws4.Range("F12").Value = 21
$ws4.Range("F15").Value = 671
$ws4.Range("F16").Value = 34
$ws4.Range("F17").Value = 19
$ws4.Range("F18").Value = 2093
$ws4.Range("F19").Value = 64
$ws4.Range("F21").Value = 73
$ws4.Range("F25").Value = 37
$ws4.Range("F27").Value = 246
$ws4.Range("F28").Value = 287
$ws4.Range("F29").Value = 419
$ws4.Range("F30").Value = 753
$ws4.Range("F33").Value = 15
